# Update "Pais" worksheet (COVID country stats) to the 22-Abril-2020 07:22 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 07:22"

# 2) In-place numeric refreshes for a handful of countries whose rank did not change.
#    (Row -> column letter -> new value)
$ws.Range("E62").Value = 1648   # Hungria: Casos activos
$ws.Range("G62").Value = 12     # Hungria: Casos criticos
$ws.Range("H62").Value = 225    # Hungria: Muertes

$ws.Range("B68").Value = 1692   # Uzbekistan: Casos totales
$ws.Range("C68").Value = 14     # Uzbekistan: Nuevos casos
$ws.Range("E68").Value = 1329   # Uzbekistan: Casos activos

$ws.Range("B75").Value = 1370   # Lituania: Casos totales
$ws.Range("C75").Value = 20     # Lituania: Nuevos casos
$ws.Range("D75").Value = 357    # Lituania: Recuperados
$ws.Range("E75").Value = 975    # Lituania: Casos activos

$ws.Range("B85").Value = 1015   # Bulgaria: Casos totales
$ws.Range("C85").Value = 40     # Bulgaria: Nuevos casos
$ws.Range("D85").Value = 174    # Bulgaria: Recuperados
$ws.Range("E85").Value = 794    # Bulgaria: Casos activos
$ws.Range("G85").Value = 2      # Bulgaria: Casos criticos
$ws.Range("H85").Value = 47     # Bulgaria: Muertes

# 3) Burundi's case count jumped from 5 to 11, moving it up the (descending) sort
#    from row 211 to row 199. Rows 199-211 shift down by one, and the former
#    "Sudan del Sur" / "Santo Tome y Principe" pair (rows 212-213) swaps order.
$rows = @(
  @{ Row=199; Country="Burundi";                      B=11; C=6; D=4;  E=6; F=0; G=0; H=1 }
  @{ Row=200; Country="Groenlandia";                   B=11; C=0; D=11; E=0; F=0; G=0; H=0 }
  @{ Row=201; Country="Gambia";                        B=10; C=0; D=2;  E=7; F=0; G=0; H=1 }
  @{ Row=202; Country="Surinam";                       B=10; C=0; D=6;  E=3; F=0; G=0; H=1 }
  @{ Row=203; Country="Nicaragua";                     B=10; C=0; D=7;  E=1; F=0; G=0; H=2 }
  @{ Row=204; Country="Santa Sede";                    B=9;  C=0; D=2;  E=7; F=0; G=0; H=0 }
  @{ Row=205; Country="Papua Nueva Guinea";             B=7;  C=0; D=0;  E=7; F=0; G=0; H=0 }
  @{ Row=206; Country="Mauritania";                    B=7;  C=0; D=6;  E=0; F=0; G=0; H=1 }
  @{ Row=207; Country="Sahara Occidental";              B=6;  C=0; D=0;  E=6; F=0; G=0; H=0 }
  @{ Row=208; Country="Butan";                         B=6;  C=0; D=2;  E=4; F=0; G=0; H=0 }
  @{ Row=209; Country="San Bartolome";                  B=6;  C=0; D=6;  E=0; F=0; G=0; H=0 }
  @{ Row=210; Country="Bonaire, San Eustaquio y Saba";  B=5;  C=0; D=0;  E=5; F=0; G=0; H=0 }
  @{ Row=211; Country="Islas Virgenes Britanicas";      B=5;  C=0; D=3;  E=1; F=0; G=0; H=1 }
  @{ Row=212; Country="Santo Tome y Principe";          B=4;  C=0; D=0;  E=4; F=0; G=0; H=0 }
  @{ Row=213; Country="Sudan del Sur";                  B=4;  C=0; D=0;  E=4; F=0; G=0; H=0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Country
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}
